$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2: BNB
Set-TextValue $ws.Range("D2") "261.43"
Set-TextValue $ws.Range("E2") "0.10%"

# Row 3: OKB
Set-TextValue $ws.Range("D3") "26.93"
Set-TextValue $ws.Range("E3") "-1.34%"

# Row 4: HuobiToken
Set-TextValue $ws.Range("D4") "4.716"
Set-TextValue $ws.Range("E4") "0.14%"

# Row 5: Cronos
Set-TextValue $ws.Range("D5") "0.06217"
Set-TextValue $ws.Range("E5") "2.26%"

# Row 6: KuCoinToken
Set-TextValue $ws.Range("D6") "6.742"
Set-TextValue $ws.Range("E6") "1.04%"

# Row 7: MXToken
Set-TextValue $ws.Range("D7") "0.8498"
Set-TextValue $ws.Range("E7") "0.47%"

# Row 8: FTXToken
Set-TextValue $ws.Range("D8") "0.9128"
Set-TextValue $ws.Range("E8") "-1.33%"

# Row 9: WazirX
Set-TextValue $ws.Range("D9") "0.1404"
Set-TextValue $ws.Range("E9") "0.04%"

# Row 10: LiechtensteinCryptoassetsExchange
Set-TextValue $ws.Range("D10") "0.04951"
Set-TextValue $ws.Range("E10") "3.90%"

# Row 11: MandalaExchangeToken
Set-TextValue $ws.Range("D11") "0.07094"
Set-TextValue $ws.Range("E11") "-0.10%"

# Row 12: BitrueCoin
Set-TextValue $ws.Range("D12") "0.03115"
Set-TextValue $ws.Range("E12") "0.65%"

# Row 14: BitForexToken
Set-TextValue $ws.Range("E14") "-0.11%"

# Row 15: One
Set-TextValue $ws.Range("D15") "0.0006178"
Set-TextValue $ws.Range("E15") "1.51%"

# Row 16: TigerCash
Set-TextValue $ws.Range("D16") "0.005964"
Set-TextValue $ws.Range("E16") "-3.54%"

# Row 17: LEO
Set-TextValue $ws.Range("D17") "3.448"
Set-TextValue $ws.Range("E17") "-0.03%"

# Row 18: GateToken
Set-TextValue $ws.Range("D18") "3.172"
Set-TextValue $ws.Range("E18") "1.08%"

# Row 19: BTSEToken
Set-TextValue $ws.Range("E19") "0.18%"

# Row 21: ProBitToken
Set-TextValue $ws.Range("D21") "0.1309"
Set-TextValue $ws.Range("E21") "1.63%"

# Row 22: MCDex
Set-TextValue $ws.Range("D22") "4.097"
Set-TextValue $ws.Range("E22") "-0.35%"

# Row 23: CoinExToken
Set-TextValue $ws.Range("D23") "0.04241"
Set-TextValue $ws.Range("E23") "0.04%"

# Row 24: BitKan
Set-TextValue $ws.Range("D24") "0.001185"
Set-TextValue $ws.Range("E24") "-3.05%"

# Row 25: HotbitToken
Set-TextValue $ws.Range("D25") "0.004073"
Set-TextValue $ws.Range("E25") "4.14%"

# Row 26: NitroEx
Set-TextValue $ws.Range("E26") "0.01%"

# Row 27: UpBots
Set-TextValue $ws.Range("E27") "4.10%"

# Row 40: IDEX
Set-TextValue $ws.Range("D40") "0.03947"
Set-TextValue $ws.Range("E40") "2.00%"

# Row 41: BKEXToken
Set-TextValue $ws.Range("D41") "0.1113"
Set-TextValue $ws.Range("E41") "0.04%"

# Row 42: KickToken
Set-TextValue $ws.Range("D42") "0.004173"
Set-TextValue $ws.Range("E42") "1.48%"

# Row 44: LocalTraders
Set-TextValue $ws.Range("D44") "0.01317"
Set-TextValue $ws.Range("E44") "-19.54%"

# Row 45: CoinLion
Set-TextValue $ws.Range("D45") "0.00005162"
Set-TextValue $ws.Range("E45") "0.12%"

# Row 46: Kangarootoken
Set-TextValue $ws.Range("E46") "0.02%"

# Row 48: BOLO
Set-TextValue $ws.Range("D48") "0.2516"
Set-TextValue $ws.Range("E48") "85.96%"

# Row 49: CryptobidCoin
Set-TextValue $ws.Range("E49") "0.02%"

# Row 50: SpecialPowerGold
Set-TextValue $ws.Range("E50") "0.02%"
